$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the paragraph "Put region name: us-east-1" by its text.
# ---------------------------------------------------------------------
$regionPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("Put region name")) {
        $regionPara = $para
    }
}

# ---------------------------------------------------------------------
# Change 1a: turn "...us-east-1" into "...us-east-1 for PSC" by
# appending two new runs just before the paragraph mark.
# ---------------------------------------------------------------------
$r = $regionPara.Range
$rTail = $d.Range($r.End - 1, $r.End - 1)
$rTail.InsertAfter(" for ")

$r2 = $regionPara.Range
$rTail2 = $d.Range($r2.End - 1, $r2.End - 1)
$rTail2.InsertAfter("PSC")

# ---------------------------------------------------------------------
# Change 1b: insert a brand-new paragraph right after it that reads
# "For Rougarou, put: us-west-2".
# ---------------------------------------------------------------------
$r3 = $regionPara.Range
$r3.InsertParagraphAfter()

$newPara = $regionPara.Next()
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.InsertAfter("For Rougarou, put: ")

$newPara2 = $regionPara.Next()
$newTail = $d.Range($newPara2.Range.End - 1, $newPara2.Range.End - 1)
$newTail.InsertAfter("us-west-2")

# ---------------------------------------------------------------------
# Change 2: merge the trailing `";  done` run with the following
# two-space run into a single run `";  done  ` (pure text, the visible
# content is unchanged - only the run split is simplified).
# ---------------------------------------------------------------------
$tailStr = '";  done  '
$tailLen = $tailStr.Length

$mergePara = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t.Contains("LWI_Production_TCs_forPurdue/`$f/") -and $t.EndsWith("done  ")) {
        $mergePara = $para
    }
}

$pr = $mergePara.Range
$rMerge = $d.Range($pr.End - 1 - $tailLen, $pr.End - 1)
$rMerge.Delete()

$pr2 = $mergePara.Range
$rInsert = $d.Range($pr2.End - 1, $pr2.End - 1)
$rInsert.InsertAfter($tailStr)
